$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data table (rows 2-11, columns B-G) after adding the latest ifoCAST
# quarter evaluation: the rolling window shifts up by one row (dropping the
# oldest quarter, Q0) and a new row of values is appended at the bottom
# (replacing Q9, now carrying the newest evaluation numbers).
$data = @(
    @(-0.05889918563118277, 0.5375210360923134, 0.6208141610631079, 0.7879176105806418, 0.8084920911945981, 18),
    @(0.07850679787661623, 0.4200275287208594, 0.2679724280602394, 0.5176605336127522, 0.5274203069366139, 17),
    @(0.09045959314081749, 0.4284335755769144, 0.3144213805930146, 0.560732895943349, 0.571536821889252, 16),
    @(0.1328729943200521, 0.4043956606847117, 0.2696880151153669, 0.5193149479028761, 0.5196490393550559, 15),
    @(0.1892581571789815, 0.3380498083070491, 0.2385894900054021, 0.4884562314121933, 0.4672992583619074, 14),
    @(0.2086948502624945, 0.3022036763864308, 0.1923865449580635, 0.4386189062934513, 0.4015417223668171, 13),
    @(0.2394124698171911, 0.2672986829739259, 0.1375220876090255, 0.3708397060847524, 0.2957955502247975, 12),
    @(0.195371401722039, 0.2718328173942094, 0.1201690274152207, 0.3466540457216975, 0.3003313954364828, 11),
    @(0.1436448583919225, 0.2781300272604909, 0.2246355840812312, 0.4739573652568669, 0.4760972574996225, 10),
    @(0.1250455391607497, 0.2940678690837308, 0.2019405311378496, 0.449377937974095, 0.4578123658313282, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $col = $c + 2
        $ws.Cells.Item($row, $col).Value = $vals[$c]
    }
}
